$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Cebollín" series.
# It becomes the new row 420, pushing every following row down by one.
$ws.Rows.Item(420).Insert()

$ws.Cells.Item(420,1).Value = 4
$ws.Cells.Item(420,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(420,3).Value = "Los Lagos"
$ws.Cells.Item(420,4).Value = 45041
$ws.Cells.Item(420,5).Value = 10
$ws.Cells.Item(420,6).Value = 100112037
$ws.Cells.Item(420,7).Value = "Cebollín"
$ws.Cells.Item(420,8).Value = "Sin especificar"
$ws.Cells.Item(420,9).Value = "Primera"
$ws.Cells.Item(420,10).Value = 160
$ws.Cells.Item(420,11).Value = 6500
$ws.Cells.Item(420,12).Value = 7000
$ws.Cells.Item(420,13).Value = 6750
$ws.Cells.Item(420,14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(420,15).Value = "Región Metropolitana"
$ws.Cells.Item(420,16).Value = 188
$ws.Cells.Item(420,17).Value = 36
$ws.Cells.Item(420,18).Value = "Hortaliza"
